$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the current row 1 data (nen, ti, 24/10/2025, 17:29, Present) down to row 2,
# update the time, and status.
$ws.Range("A2").Value = "nen"
$ws.Range("B2").Value = "ti"
$ws.Range("C2").Value = "24/10/2025"
$ws.Range("D2").Value = "17:43"
$ws.Range("E2").Value = "Absent"

# Make sure the new data row (row 2) has no special styling (default style).
$ws.Range("A2:E2").Style = "Normal"

# Set the new header row (row 1) values.
$ws.Range("A1").Value = "Student"
$ws.Range("B1").Value = "Course"
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "Time"
$ws.Range("E1").Value = "Status"

# Re-apply the header style to row 1 (it already had style index 1 applied: bold, bordered, centered).
$ws.Range("A1:E1").Font.Bold = $true
$ws.Range("A1:E1").Borders.LineStyle = 1
$ws.Range("A1:E1").HorizontalAlignment = -4108
$ws.Range("A1:E1").VerticalAlignment = -4160
